# This script applies a weekly update to the "Berenjena" price sheet:
# a new week's record is inserted as row 108 (pushing the existing
# rows 108-195 down to 109-196), and the new row is populated with
# the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 108; everything from the old
# row 108 downward shifts down by one row.
$ws.Rows.Item(108).Insert()

# Populate the new row 108 with the new week's record. Most fields
# are identical to the (now shifted) row 109, only the date, volume,
# min/max price and weighted average price (and the derived $/Kg)
# differ for the new week.
$ws.Range("A108").Value = 6
$ws.Range("B108").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C108").Value = "Metropolitana"
$ws.Range("D108").Value = 44651
$ws.Range("D108").NumberFormat = $ws.Range("D109").NumberFormat
$ws.Range("E108").Value = 13
$ws.Range("F108").Value = 100112001
$ws.Range("G108").Value = "Berenjena"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 430
$ws.Range("K108").Value = 8000
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 8395
$ws.Range("N108").Value = "$/caja 50 unidades"
$ws.Range("O108").Value = "Región de Arica y Parinacota"
$ws.Range("P108").Value = 168
$ws.Range("Q108").Value = 50
$ws.Range("R108").Value = "Hortaliza"
